# Kama.Calc.xlsx refactor: rename the "index" column/header to "i" and
# shift the running index values from 1-based to 0-based (rows 2..503).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header "index" -> "i" (column A, also the first column of the
# "testdata" table, whose column name tracks the header cell's text).
$ws.Range("A1").Value = "i"

# The data rows used to start the running counter at 1 (row 2 -> 1, row 3 -> 2,
# ... row 503 -> 502). Shift every value down by one so it starts at 0.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() - 1
}

# The column got visually narrower once the header/content shortened; match
# the resulting best-fit width.
$ws.Columns.Item(1).ColumnWidth = 3.14
